$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 317; this shifts the existing rows 317-338 down to 318-339,
# carrying all their original data (and formatting) with them.
$ws.Rows.Item(317).Insert()

# Populate the newly inserted row 317 with the new weekly price record.
$ws.Cells.Item(317, 1).Value = 4
$ws.Cells.Item(317, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(317, 3).Value = "Los Lagos"
$ws.Cells.Item(317, 4).Value = 44826
$ws.Cells.Item(317, 5).Value = 10
$ws.Cells.Item(317, 6).Value = 100112037
$ws.Cells.Item(317, 7).Value = "Cebollín"
$ws.Cells.Item(317, 8).Value = "Sin especificar"
$ws.Cells.Item(317, 9).Value = "Primera"
$ws.Cells.Item(317, 10).Value = 100
$ws.Cells.Item(317, 11).Value = 8000
$ws.Cells.Item(317, 12).Value = 9000
$ws.Cells.Item(317, 13).Value = 8500
$ws.Cells.Item(317, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(317, 15).Value = "Región Metropolitana"
$ws.Cells.Item(317, 16).Value = 236
$ws.Cells.Item(317, 17).Value = 36
$ws.Cells.Item(317, 18).Value = "Hortaliza"
